$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 22.59487733333333
$ws.Range("H2").Value = 67.784632
$ws.Range("I2").Value = 0.7395019553569895
$ws.Range("J2").Value = 0.7395019553569895
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.08532
$ws.Range("N2").Value = 6.25596
$ws.Range("O2").Value = 0.01753772176136817
$ws.Range("P2").Value = 0.01753772176136816
$ws.Range("Q2").Value = 47.11754960074666
$ws.Range("R2").Value = 424.05794640672
$ws.Range("S2").Value = 0.01296917953503858
$ws.Range("T2").Value = 0.01296917953503858

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 22.59487733333333
$ws.Range("H3").Value = 67.784632
$ws.Range("I3").Value = 0.7395019553569895
$ws.Range("J3").Value = 0.7395019553569895
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.898173
$ws.Range("N3").Value = 305.694519
$ws.Range("O3").Value = 0.8569724579756384
$ws.Range("P3").Value = 0.8569724579756383
$ws.Range("Q3").Value = 2302.376719425778
$ws.Range("R3").Value = 20721.39047483201
$ws.Range("S3").Value = 0.6337328083600701
$ws.Range("T3").Value = 0.63373280836007

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 22.59487733333333
$ws.Range("H4").Value = 67.784632
$ws.Range("I4").Value = 0.7395019553569895
$ws.Range("J4").Value = 0.7395019553569895
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.921347
$ws.Range("N4").Value = 44.76404100000001
$ws.Range("O4").Value = 0.1254898202629935
$ws.Range("P4").Value = 0.1254898202629935
$ws.Range("Q4").Value = 337.1460051131014
$ws.Range("R4").Value = 3034.314046017912
$ws.Range("S4").Value = 0.09279996746188086
$ws.Range("T4").Value = 0.09279996746188085

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3045986666666667
$ws.Range("H5").Value = 0.9137960000000001
$ws.Range("I5").Value = 0.00996913177602551
$ws.Range("J5").Value = 0.00996913177602551
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.08532
$ws.Range("N5").Value = 6.25596
$ws.Range("O5").Value = 0.01753772176136817
$ws.Range("P5").Value = 0.01753772176136816
$ws.Range("Q5").Value = 0.6351856915733333
$ws.Range("R5").Value = 5.716671224160001
$ws.Range("S5").Value = 0.0001748358592903495
$ws.Range("T5").Value = 0.0001748358592903494

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.3045986666666667
$ws.Range("H6").Value = 0.9137960000000001
$ws.Range("I6").Value = 0.00996913177602551
$ws.Range("J6").Value = 0.00996913177602551
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.898173
$ws.Range("N6").Value = 305.694519
$ws.Range("O6").Value = 0.8569724579756384
$ws.Range("P6").Value = 0.8569724579756383
$ws.Range("Q6").Value = 31.03804763156933
$ws.Range("R6").Value = 279.342428684124
$ws.Range("S6").Value = 0.008543271361983623
$ws.Range("T6").Value = 0.008543271361983622

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.3045986666666667
$ws.Range("H7").Value = 0.9137960000000001
$ws.Range("I7").Value = 0.00996913177602551
$ws.Range("J7").Value = 0.00996913177602551
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.921347
$ws.Range("N7").Value = 44.76404100000001
$ws.Range("O7").Value = 0.1254898202629935
$ws.Range("P7").Value = 0.1254898202629935
$ws.Range("Q7").Value = 4.545022401070668
$ws.Range("R7").Value = 40.90520160963601
$ws.Range("S7").Value = 0.001251024554751539
$ws.Range("T7").Value = 0.001251024554751538

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.654706
$ws.Range("H8").Value = 22.964118
$ws.Range("I8").Value = 0.2505289128669849
$ws.Range("J8").Value = 0.2505289128669849
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.08532
$ws.Range("N8").Value = 6.25596
$ws.Range("O8").Value = 0.01753772176136817
$ws.Range("P8").Value = 0.01753772176136816
$ws.Range("Q8").Value = 15.96251151592
$ws.Range("R8").Value = 143.66260364328
$ws.Range("S8").Value = 0.00439370636703923
$ws.Range("T8").Value = 0.00439370636703923

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.654706
$ws.Range("H9").Value = 22.964118
$ws.Range("I9").Value = 0.2505289128669849
$ws.Range("J9").Value = 0.2505289128669849
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.898173
$ws.Range("N9").Value = 305.694519
$ws.Range("O9").Value = 0.8569724579756384
$ws.Range("P9").Value = 0.8569724579756383
$ws.Range("Q9").Value = 780.000556252138
$ws.Range("R9").Value = 7020.005006269242
$ws.Range("S9").Value = 0.2146963782535846
$ws.Range("T9").Value = 0.2146963782535846

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.654706
$ws.Range("H10").Value = 22.964118
$ws.Range("I10").Value = 0.2505289128669849
$ws.Range("J10").Value = 0.2505289128669849
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.921347
$ws.Range("N10").Value = 44.76404100000001
$ws.Range("O10").Value = 0.1254898202629935
$ws.Range("P10").Value = 0.1254898202629935
$ws.Range("Q10").Value = 114.218524408982
$ws.Range("R10").Value = 1027.966719680838
$ws.Range("S10").Value = 0.0314388282463611
$ws.Range("T10").Value = 0.03143882824636109

